$wb = $excel.ActiveWorkbook

# --- Sheet 1: Sh_exact ---
$ws1 = $wb.Worksheets.Item("Sh_exact")
$ws1.Range("B2").Value = 0.104813269303049
$ws1.Range("C2").Value = 0.09585042201537033
$ws1.Range("D2").Value = 0.1092849281569675

$ws1.Range("B3").Value = 0.4171261421903749
$ws1.Range("C3").Value = 0.4149373803962081
$ws1.Range("D3").Value = 0.4209267934553773

$ws1.Range("B4").Value = 0.4780605885065762
$ws1.Range("C4").Value = 0.4751933081685998
$ws1.Range("D4").Value = 0.483669597729377

# --- Sheet 2: Vsob_exact ---
$ws2 = $wb.Worksheets.Item("Vsob_exact")
$ws2.Range("B2").Value = 0.1049220128443874
$ws2.Range("C2").Value = 0.07138096436810831
$ws2.Range("D2").Value = 0.1214878237641213

$ws2.Range("B3").Value = 0.8155801991930413
$ws2.Range("C3").Value = 0.8088032668573277
$ws2.Range("D3").Value = 0.8189034807945491

$ws2.Range("B4").Value = 0.8762204971247777
$ws2.Range("C4").Value = 0.8714746624616878
$ws2.Range("D4").Value = 0.878456829485212

# --- Sheet 3: Tsob_exact ---
$ws3 = $wb.Worksheets.Item("Tsob_exact")
$ws3.Range("B2").Value = 0.1050581125744231
$ws3.Range("C2").Value = 0.1031826900439379
$ws3.Range("D2").Value = 0.1091002248472033

$ws3.Range("B3").Value = 0.01973773135800618
$ws3.Range("C3").Value = 0.01936216511791973
$ws3.Range("D3").Value = 0.02048436930152382

$ws3.Range("B4").Value = 0.0793823750306652
$ws3.Range("C4").Value = 0.07788345051162671
$ws3.Range("D4").Value = 0.08247574127698191
